$wb = $excel.ActiveWorkbook

# Update "想去人数" (F2, F5) on both "展览" and "全部类型" sheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 44
    $ws.Range("F5").Value = 51
}
